$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9, shifting rows 9 and below down by one.
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the "Deposito" entry.
$ws.Range("B9").Value = "Deposito"
$ws.Range("C9").Value = "idDeposito"
$ws.Range("E9").Value = "idFormaDePago → FormaDePago"

# Update selection to match target workbook state.
$ws.Range("E9").Select()
